# Appends the new "Exploration / campaign mode" section to the end of
# the design document, matching the authored diff:
#   - two new Heading1 paragraphs ("Exploration / campaign mode",
#     "On implementation")
#   - five new body paragraphs discussing the implementation notes
#
# The new content is inserted as plain runs inheriting the surrounding
# en-US language formatting (matching the rest of the document), and the
# two heading paragraphs are then promoted to the "Heading1" style.

$d = $word.ActiveDocument

# Anchor a collapsed range at the very end of the main document story so
# the new text lands after the last existing paragraph ("Each square has
# a pressure variable...") and before the trailing sectPr.
$endRange = $d.Range($d.Content.End, $d.Content.End)

# Leading `r makes sure the first new paragraph does not get appended
# onto the existing last paragraph's run.
$newText = "`rExploration / campaign mode"
$newText += "`rOn implementation"
$newText += "`rThe graphics are only 2D to eliminate one dimension of complexity. The best way is probably to go with OpenGL (+SDL?) for the ease of programming and multiplatform. No pre-made engines unless required for maximum learning experience. Physics should prove to be a fun puzzle, and the design might change depending on the outcome of the battle with the physics."
$newText += "`rThe battle arena is probably filled with empty space. There are some 1000 - 10 000 squares inside the ship, and inside each square there are components, humans and other things of note. This means that there will be quite a load of variables changing for each frame, so perhaps some kind of list for `"updated objects`" is required. Graphics will probably be simplified first with the ship showing the floor plan with nothing moving inside. The humans and other stuff can be added later."
$newText += "`rThe squares inside the ship could probably be defined relative to the hull, not screen borders. Each square could either have the X and Y coordinates on the screen, or they might simply be calculated by a simple algorithm, or they might be in a 2-dimensional array."
$newText += " Anyway, they need to be set in such fashion that they are always static in relation to the ship, but so that the squares can change their appearance (damage makes them black or something)."
$newText += " It might be "
$newText += "reasonable to build the ship from a set of flat 2D planes with a texture on them as this way rendering the geometry is fast. In future, the third dimension could be added for the look (2.5D), and the ships could simply move in XY-coordinate system."
$newText += "`rFor the starters it might be good to ignore any thoughts on the campaign / exploration mode. First there should simply be the battle mode where the ship can move around, then maybe try some collision tests, then try getting some of the module code in."

$endRange.InsertAfter($newText)

# The two heading paragraphs were inserted as plain paragraphs; give them
# the same "Heading1" style used for the other top-level section titles
# ("Battle mode", "About this document", ...). We locate them by their
# text instead of a hard-coded index so the script is resilient.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Exploration / campaign mode`r" -or $t -eq "On implementation`r") {
        $p.Style = "Heading1"
    }
}

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)

